$d = $word.ActiveDocument

# --- Change 1: remove the duplicate 4-space run right before "demonstration" ---
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Replacement.ClearFormatting()
[void]$rng1.Find.Execute("    demonstration", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "demonstration", 2)

# --- Change 2: turn the `{ m:self.name }` field into plain literal text runs ---
# Locate the field (it is the only field in the document) and find which
# paragraph holds it, so we can come back to that same paragraph once the
# field itself has been deleted.
$field = $d.Fields.Item(1)
$fieldStart = $field.Code.Start
$paraIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($fieldStart -ge $p.Range.Start -and $fieldStart -le $p.Range.End) {
        $paraIndex = $i
        break
    }
}

# Delete the field (removes the fldChar/instrText runs entirely).
$field.Delete()

# Re-fetch the (now empty) paragraph that used to hold the field.
$para = $d.Paragraphs.Item($paraIndex)

# Insert the replacement runs (plain <w:t> runs, "self" keeps its orange theme color)
# directly as OOXML so the exact theme-color markup is preserved.
$xmlFrag = '<?xml version="1.0" standalone="yes"?>' + `
  '<?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' + `
      '<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' + `
        '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' + `
      '</Relationships></pkg:xmlData>' + `
    '</pkg:part>' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + `
          '<w:r><w:t>{</w:t></w:r>' + `
          '<w:r><w:t>m</w:t></w:r>' + `
          '<w:r><w:t>:</w:t></w:r>' + `
          '<w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>self</w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve">.name}</w:t></w:r>' + `
        '</w:p></w:body></w:document></pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

[void]$para.Range.InsertXML($xmlFrag)
